# Swap the order of the two comma-separated "Recorded By" values in
# column G: "dnasr281@gmail.com, X" -> "X, dnasr281@gmail.com"
# (affects rows where the first recorder listed is dnasr281@gmail.com
#  and there is exactly one other recorder, e.g. "System" or
#  "admin@admin.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $txt = $cell.Text

    if ($txt -like "dnasr281@gmail.com, *") {
        $parts = $txt -split ", "
        if ($parts.Count -eq 2) {
            $swapped = $parts[1] + ", " + $parts[0]
            $cell.Value = $swapped
        }
    }
}
